# Updates cryptocurrency price/volume table (Coin, Link, Price, Volume(1h))
# with freshly scraped values, matching GitHub Actions scheduled refresh.
# Rows 36-38 also get re-ranked (InternetComputer/VeChain/Algorand swap order).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.103.38"
$ws.Range("E2").Value = "  -3.31%  "
$ws.Range("D3").Value = "1.599.53"
$ws.Range("E3").Value = "  -3.04%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.001"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "301.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.00%  "
$ws.Range("E7").Value = "  -2.81%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3641"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "49.88"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.262"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.001"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08138"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.78%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.02"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.578"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001257"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.359"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -8.25%  "
$ws.Range("D17").Value = "1.601.70"
$ws.Range("E17").Value = "  -3.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06861"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.552"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.73%  "
$ws.Range("D24").Value = "23.097.14"
$ws.Range("E24").Value = "  -3.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.339"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.712"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.273"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.03"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.415"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.819"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -12.80%  "
$ws.Range("D33").Value = "1.778.36"
$ws.Range("E33").Value = "  -2.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9565"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.94%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07639"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.70%  "
$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.249"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.88%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02722"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.27%  "
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2543"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.08896"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.06"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.367"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7079"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.64"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6590"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.08%  "
$ws.Range("E46").Value = "  +0.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.303"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.981"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "130.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07936"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.207"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.04%  "

Write-Output "Updated crypto values for $([DateTime]::Now)"
